# "vault backup: 2023-01-04 14:56:55"
# Row 3 (Puerto Limon replaces "Meet the sloths"), cheaper single/double
# price, and a new "Select" (D) column value; the grand-total row gets
# recomputed and its label is rewritten as a plain cell value; the sheet
# view/column widths are nudged too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: new excursion name + lower single/double prices ---------------
$ws.Range("A3").Value = "Puerto Limon, Costa Rica (Cultural & Soul food)"
$ws.Range("B3").Value = 90
# C3 already holds "=B3*2" and recalculates on its own.

# D3 is a new "Select" cell (mirrors D4, D5, D7, D8 = "=C<row>"); copy the
# existing D4 formatting first so it reuses the same cell style rather than
# minting a new one, then drop in the formula.
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("D3").Formula = "=C3"

# --- Grand total row (27): label becomes a literal "Total" value ----------
$ws.Range("A27").Value = "Total"
# B27/C27/D27 keep their existing SUM formulas; they recompute automatically
# once B3/D3 change.

# --- Column widths: A widens, D widens (no longer "best fit") -------------
$ws.Columns.Item(1).ColumnWidth = 43.30687099358975
$ws.Columns.Item(4).ColumnWidth = 15.73709935897436

# --- Sheet view: drop the frozen top-left cell, move the selection --------
$ws.Range("B4").Select()

$excel.CalculateFull()
